# Auto-generated edit script applying scheduled-runner value updates
# to the Ultros_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 732.8333
$ws.Range("I2").Value = 824
$ws.Range("J2").Value = 687.25
$ws.Range("K2").Value = 824
$ws.Range("L2").Value = 687.25
$ws.Range("M2").Value = -711
$ws.Range("N2").Value = -913.25
$ws.Range("H9").Value = 155.6
$ws.Range("I9").Value = 155.6
$ws.Range("K9").Value = 155.6
$ws.Range("M9").Value = 13.40000000000001
$ws.Range("H87").Value = 49999.934
$ws.Range("J87").Value = 49999.934
$ws.Range("L87").Value = 49999.934
$ws.Range("N87").Value = -52495.934
$ws.Range("H90").Value = 49999.934
$ws.Range("J90").Value = 49999.934
$ws.Range("L90").Value = 149999.802
$ws.Range("N90").Value = -162479.802
$ws.Range("H113").Value = 6973.3145
$ws.Range("I113").Value = 6190.2
$ws.Range("J113").Value = 8017.467
$ws.Range("K113").Value = 6190.2
$ws.Range("L113").Value = 8017.467
$ws.Range("M113").Value = -2936.2
$ws.Range("N113").Value = -14525.467
$ws.Range("H125").Value = 3404.2222
$ws.Range("I125").Value = 1434.5
$ws.Range("J125").Value = 4980
$ws.Range("K125").Value = 12910.5
$ws.Range("L125").Value = 44820
$ws.Range("M125").Value = -10450.5
$ws.Range("N125").Value = -49740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21165.533
$ws.Range("I2").Value = 31245.1
$ws.Range("K2").Value = 31245.1
$ws.Range("M2").Value = -31132.1
$ws.Range("H32").Value = 23074978
$ws.Range("I32").Value = 25021016
$ws.Range("K32").Value = 25021016
$ws.Range("M32").Value = -25020729
$ws.Range("H61").Value = 4969.7856
$ws.Range("I61").Value = 2956.2942
$ws.Range("K61").Value = 2956.2942
$ws.Range("M61").Value = -2744.2942
$ws.Range("H74").Value = 5168.7
$ws.Range("I74").Value = 5316.3335
$ws.Range("J74").Value = 4947.25
$ws.Range("K74").Value = 5316.3335
$ws.Range("L74").Value = 4947.25
$ws.Range("M74").Value = -4442.3335
$ws.Range("N74").Value = -6695.25
$ws.Range("H77").Value = 5168.7
$ws.Range("I77").Value = 5316.3335
$ws.Range("J77").Value = 4947.25
$ws.Range("K77").Value = 26581.6675
$ws.Range("L77").Value = 24736.25
$ws.Range("M77").Value = -22213.6675
$ws.Range("N77").Value = -33472.25
$ws.Range("H102").Value = 2590.3333
$ws.Range("I102").Value = 2108.4
$ws.Range("K102").Value = 2108.4
$ws.Range("M102").Value = -486.4000000000001
$ws.Range("H116").Value = 21165.533
$ws.Range("I116").Value = 31245.1
$ws.Range("K116").Value = 31245.1
$ws.Range("M116").Value = -28951.1
$ws.Range("H132").Value = 4343.1978
$ws.Range("I132").Value = 3663.288
$ws.Range("J132").Value = 5828.926
$ws.Range("K132").Value = 10989.864
$ws.Range("L132").Value = 17486.778
$ws.Range("M132").Value = -8459.864
$ws.Range("N132").Value = -22546.778
$ws.Range("H136").Value = 4969.7856
$ws.Range("I136").Value = 2956.2942
$ws.Range("K136").Value = 8868.882599999999
$ws.Range("M136").Value = -6318.882599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21165.533
$ws.Range("I3").Value = 31245.1
$ws.Range("K3").Value = 31245.1
$ws.Range("M3").Value = -31131.1
$ws.Range("H99").Value = 22878.264
$ws.Range("I99").Value = 41235.6
$ws.Range("K99").Value = 41235.6
$ws.Range("M99").Value = -39737.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3928.2727
$ws.Range("J58").Value = 6750
$ws.Range("L58").Value = 6750
$ws.Range("N58").Value = -7156
$ws.Range("H99").Value = 7089.8096
$ws.Range("I99").Value = 7066.2
$ws.Range("K99").Value = 7066.2
$ws.Range("M99").Value = -5568.2
$ws.Range("H126").Value = 7089.8096
$ws.Range("I126").Value = 7066.2
$ws.Range("K126").Value = 21198.6
$ws.Range("M126").Value = -18728.6
$ws.Range("H132").Value = 4911.0312
$ws.Range("I132").Value = 4367.1875
$ws.Range("K132").Value = 13101.5625
$ws.Range("M132").Value = -10571.5625
$ws.Range("H134").Value = 5941.4
$ws.Range("I134").Value = 4147.231
$ws.Range("K134").Value = 12441.693
$ws.Range("M134").Value = -9906.692999999999
$ws.Range("H136").Value = 3928.2727
$ws.Range("J136").Value = 6750
$ws.Range("L136").Value = 20250
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3401.4
$ws.Range("J113").Value = 3373.625
$ws.Range("L113").Value = 10120.875
$ws.Range("N113").Value = -14460.875
$ws.Range("H131").Value = 5216.304
$ws.Range("J131").Value = 5441.385
$ws.Range("L131").Value = 16324.155
$ws.Range("N131").Value = -26404.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 121857.14
$ws.Range("J105").Value = 121857.14
$ws.Range("L105").Value = 121857.14
$ws.Range("N105").Value = -128845.14
$ws.Range("H122").Value = 4126.2354
$ws.Range("I122").Value = 3819.6924
$ws.Range("K122").Value = 11459.0772
$ws.Range("M122").Value = -9009.0772
$ws.Range("H126").Value = 5131.4165
$ws.Range("I126").Value = 3299.5
$ws.Range("J126").Value = 5497.8
$ws.Range("K126").Value = 9898.5
$ws.Range("L126").Value = 16493.4
$ws.Range("M126").Value = -7428.5
$ws.Range("N126").Value = -21433.4
$ws.Range("H132").Value = 9040.023999999999
$ws.Range("I132").Value = 7580.759
$ws.Range("J132").Value = 12566.583
$ws.Range("K132").Value = 22742.277
$ws.Range("L132").Value = 37699.749
$ws.Range("M132").Value = -20212.277
$ws.Range("N132").Value = -42759.749

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5635.636
$ws.Range("J7").Value = 6698.625
$ws.Range("L7").Value = 6698.625
$ws.Range("N7").Value = -6922.625
$ws.Range("H40").Value = 15438.272
$ws.Range("I40").Value = 28402.75
$ws.Range("K40").Value = 28402.75
$ws.Range("M40").Value = -28266.75
$ws.Range("H42").Value = 22250
$ws.Range("J42").Value = 22250
$ws.Range("L42").Value = 22250
$ws.Range("N42").Value = -23376
$ws.Range("H49").Value = 22250
$ws.Range("J49").Value = 22250
$ws.Range("L49").Value = 22250
$ws.Range("N49").Value = -22544
$ws.Range("H61").Value = 3015.375
$ws.Range("I61").Value = 2337.8572
$ws.Range("K61").Value = 2337.8572
$ws.Range("M61").Value = -2135.8572
$ws.Range("H113").Value = 3015.375
$ws.Range("I113").Value = 2337.8572
$ws.Range("K113").Value = 2337.8572
$ws.Range("M113").Value = -167.8571999999999
$ws.Range("H122").Value = 6598.0586
$ws.Range("J122").Value = 9044.888999999999
$ws.Range("L122").Value = 27134.667
$ws.Range("N122").Value = -32034.667
$ws.Range("H126").Value = 5635.636
$ws.Range("J126").Value = 6698.625
$ws.Range("L126").Value = 20095.875
$ws.Range("N126").Value = -25035.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1881.5385
$ws.Range("I107").Value = 1834
$ws.Range("K107").Value = 5502
$ws.Range("M107").Value = -3582
$ws.Range("H122").Value = 4193.857
$ws.Range("I122").Value = 3960.6155
$ws.Range("J122").Value = 4572.875
$ws.Range("K122").Value = 11881.8465
$ws.Range("L122").Value = 13718.625
$ws.Range("M122").Value = -9431.8465
$ws.Range("N122").Value = -18618.625
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12530
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3804
$ws.Range("I132").Value = 2593.675
$ws.Range("J132").Value = 6352.0527
$ws.Range("K132").Value = 7781.025000000001
$ws.Range("L132").Value = 19056.1581
$ws.Range("M132").Value = -5251.025000000001
$ws.Range("N132").Value = -24116.1581
